$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("I2").Value = -0.5298692899008431
$ws.Range("J2").Value = 0.121416946528233
$ws.Range("K2").Value = -0.3367949059522933
$ws.Range("L2").Value = 2.03617404199618

# Row 12
$ws.Range("I12").Value = -0.2569384304019269
$ws.Range("J12").Value = 0.04487644048399081
$ws.Range("K12").Value = -0.6591431417536935
$ws.Range("L12").Value = 1.5

# Row 13
$ws.Range("I13").Value = -0.7095188395570421
$ws.Range("J13").Value = 0.1611014700843928
$ws.Range("K13").Value = 0.2599005274089113
$ws.Range("L13").Value = 1.583864174784634

# Row 14
$ws.Range("I14").Value = -0.4066391886530542
$ws.Range("J14").Value = 0.1008882814258805
$ws.Range("K14").Value = -0.3031090771618764
$ws.Range("L14").Value = 1.426911848121309

# Row 15
$ws.Range("I15").Value = -0.4372884950901768
$ws.Range("J15").Value = 0.08446814508101096
$ws.Range("K15").Value = -0.02830387734468846
$ws.Range("L15").Value = 1.565151079672559

# Row 16
$ws.Range("I16").Value = -0.4292028699684391
$ws.Range("J16").Value = 0.07071700756888885
$ws.Range("K16").Value = -0.1210231726677569
$ws.Range("L16").Value = 1.870213651501798

# Row 17
$ws.Range("I17").Value = -0.7069195875716421
$ws.Range("J17").Value = 0.1189669907999237
$ws.Range("K17").Value = 0.2629798367475658
$ws.Range("L17").Value = 2.073549511692321
